{"js": "// Rename the \"Case QC Passed\" column header to \"Treated QC Passed\", and\n// replace every standalone \"Pass\" result cell with \"yes\" (naming-convention\n// update per the commit message). \"Control QC Passed\" header is left as-is.\n\nconst body = context.document.body;\n\n// 1) Header cell: \"Case QC Passed\" -> \"Treated QC Passed\"\nconst headerHits = body.search(\"Case QC Passed\", { matchCase: true });\nheaderHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of headerHits.items) {\n  hit.insertText(\"Treated QC Passed\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Every QC result cell that reads exactly \"Pass\" -> \"yes\"\n//    matchWholeWord keeps this from touching \"...Passed\" header text.\nconst passHits = body.search(\"Pass\", { matchCase: true, matchWholeWord: true });\npassHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of passHits.items) {\n  hit.insertText(\"yes\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Rename the \"Case QC Passed\" column header to \"Treated QC Passed\", and\n# replace every standalone \"Pass\" result cell with \"yes\" (naming-convention\n# update per the commit message). \"Control QC Passed\" header is left as-is.\n\n$d = $word.ActiveDocument\n\n# 1) Header cell: \"Case QC Passed\" -> \"Treated QC Passed\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Case QC Passed\"\n$find.Replacement.Text = \"Treated QC Passed\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, \"Treated QC Passed\", 2)\n\n# 2) Every QC result cell that reads exactly \"Pass\" -> \"yes\"\n#    MatchWholeWord keeps this from touching \"...Passed\" header text.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Pass\"\n$find2.Replacement.Text = \"yes\"\n$find2.Execute($null, $true, $true, $false, $false, $false, $true, 1, $false, \"yes\", 2)\n\nWrite-Output \"done\"\n"}
